$d = $word.ActiveDocument

# Locate the paragraph that holds the "..commit 3 ..." text (old commit note)
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text -like "*commit*") {
        $target = $para
    }
}

$full = $target.Range
$body = $d.Range($full.Start, $full.End - 1)
$body.Text = "commit 4"

# White-on-white text (both the run and the paragraph mark formatting)
$body.Font.Color = 16777215
$target.Range.Font.Color = 16777215
